$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 17:34"

# Row 4
$ws.Range("B4").Value = 1296044
$ws.Range("C4").Value = 3421
$ws.Range("E4").Value = 1001591
$ws.Range("G4").Value = 233
$ws.Range("H4").Value = 77161

# Row 11
$ws.Range("B11").Value = 137309
$ws.Range("C11").Value = 1616
$ws.Range("E11").Value = 72652
$ws.Range("G11").Value = 119
$ws.Range("H11").Value = 9307

# Row 48
$ws.Range("B48").Value = 8055
$ws.Range("C48").Value = 21
$ws.Range("E48").Value = 7805
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 218

# Row 99
$ws.Range("B99").Value = 891
$ws.Range("C99").Value = 2
$ws.Range("E99").Value = 476
$ws.Range("F99").Value = 10

# Row 173
$ws.Range("B173").Value = 47
$ws.Range("C173").Value = 2
$ws.Range("E173").Value = 17

# Row 204
$ws.Range("D204").Value = 8
$ws.Range("E204").Value = 3
